$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.636.53"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "3.152.15"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +14.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "3.700.06"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "58.713.66"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.25"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "3.154.03"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.81"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.520"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +13.14%  "
$ws.Range("D29").Value = "0.0₃0862"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.13"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  +5.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").Value = "2.644.29"
$ws.Range("E40").Value = "  +6.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0686"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("E42").Value = "  +6.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.35%  "
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  +6.00%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "3.197.02"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("E48").Value = "  +13.72%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.980"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.78%  "
